$wb = $excel.ActiveWorkbook

# Rename the "Dynamic" sheet to reflect the new bank range 0x13 - 0x28
$dyn = $wb.Worksheets.Item("Dynamic (Bank 0x13 - 0x3A)")
$dyn.Name = "Dynamic (Bank 0x13 - 0x28)"

# Correct the bank allocation ranges and counts on the Dynamic sheet
$dyn.Range("E4").Value = "0x15-0x1A"
$dyn.Range("E5").Value = "0x1B-0x22"
$dyn.Range("E6").Value = "0x23-0x28"

$dyn.Range("F4").Value = 6
$dyn.Range("F5").Value = 8

# Make the Dynamic sheet the active sheet/tab, then update its selection
$dyn.Activate()
$dyn.Range("J6").Select()
